# "Volviendo light los resumenes": shorten the Spanish and English summary
# paragraphs on the back-cover (contraportada) slide of the CTFG deck.
#
# The slide's only shape holds a single paragraph shaped like:
#   run (es-ES) -> br -> br -> run (en-US) -> run (en-US) -> run (en-US)
#
# Target layout after the edit:
#   run (es-ES) -> run (es-ES) -> run (es-ES) -> br -> br -> run (en-US)
#
# i.e. the Spanish text gets trimmed and re-split into three runs around a
# newly inserted "de datos, que permitirá," clause, while the three English
# runs collapse into a single, shorter run.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$tf = $s.Shapes.Item(1).TextFrame
$tr = $tf.TextRange

$originalText = $tr.Text

# ---------------------------------------------------------------------
# Locate the English block (after the two manual line breaks) using a
# unique marker from its start, and work out its length from the known
# end-of-paragraph marker.
# ---------------------------------------------------------------------
$englishMarker = "A very important aspect"
$englishStart0 = $originalText.IndexOf($englishMarker)   # 0-based
$englishStart  = $englishStart0 + 1                      # 1-based for Characters()
$englishLen    = $originalText.Length - $englishStart0

$newEnglish = "A very important aspect that must be considered in web development is the security of the interactions that originate between the client and the server. Being able to provide the tool with security in its transactions and communication with the user is an ideal complementary aspect of the web development carried out in the final degree project. This complement aims to cover essential fields in the security of the website, such as the tokenization of queries or the protection of access to certain routes of the page. The functionalities offered by the website will also be expanded, implementing an email manager that will help in all interactions with the user and a system capable of dumping data from an Excel file. A backup system will be developed for the database that will allow the data stored in UAL Inventarium to be recovered with the minimum loss of information."

$enRange = $tr.Characters($englishStart, $englishLen)
$enRange.Text = $newEnglish

# ---------------------------------------------------------------------
# Spanish block: originally a single run, starting at character 1 and
# running up to (but excluding) the first manual line break.
# ---------------------------------------------------------------------
$breakIndex0  = $originalText.IndexOf([char]11)  # vertical-tab == <a:br/>
$spanishStart = 1
$spanishLen   = $breakIndex0                     # chars before the first break

$esPart1 = "Un aspecto muy importante que hay que considerar dentro del desarrollo web es el de la seguridad de las interacciones que se originan entre el cliente y el servidor. El poder dotar a la herramienta, de seguridad en sus transacciones y de una comunicación con el usuario, es un aspecto complementario ideal para el desarrollo web realizado en el trabajo fin de grado. Este complemento pretende cubrir campos esenciales en la seguridad del sitio web como puede ser la tokenización de las consultas o la protección del acceso a determinadas rutas de la página.  También se ampliarán las funcionalidades que ofrece la web, implementando un gestor de correos que ayudara en todas las interacciones con el usuario y un sistema capaz de realizar un volcado de datos desde un archivo Excel. Se desarrollará un sistema de copias de seguridad para la base "
$esPart2 = "de datos, que permitirá, "
$esPart3 = "poder recuperar los datos almacenados en UAL Inventarium con la menor pérdida de información. "

$newSpanish = $esPart1 + $esPart2 + $esPart3

$esRange = $tr.Characters($spanishStart, $spanishLen)
$esRange.Text = $newSpanish

# Re-touch the middle chunk ("de datos, que permitirá, ") in place so the
# paragraph ends up with three Spanish runs instead of one, matching the
# target run layout (a separate run for the newly inserted clause).
$midStart = $spanishStart + $esPart1.Length
$midLen   = $esPart2.Length
$midRange = $tr.Characters($midStart, $midLen)
$midRange.Text = $esPart2
